$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue 2 4 '26.426.22'
Set-TextValue 2 5 '  -0.38%  '
Set-TextValue 3 4 '1.724.33'
Set-TextValue 3 5 '  -0.24%  '
Set-TextValue 4 4 '0.9993'
Set-TextValue 4 5 '  +0.00%  '
Set-TextValue 5 4 '242.87'
Set-TextValue 5 5 '  -0.80%  '
Set-TextValue 6 4 '0.9996'
Set-TextValue 6 5 '  -0.01%  '
Set-TextValue 7 4 '0.4916'
Set-TextValue 7 5 '  +2.19%  '
Set-TextValue 8 4 '0.2618'
Set-TextValue 8 5 '  -2.10%  '
Set-TextValue 9 4 '0.06205'
Set-TextValue 9 5 '  +0.36%  '
Set-TextValue 10 4 '1.712.32'
Set-TextValue 10 5 '  -0.85%  '
Set-TextValue 11 4 '0.07023'
Set-TextValue 11 5 '  -2.22%  '
Set-TextValue 12 4 '15.50'
Set-TextValue 12 5 '  -0.53%  '
Set-TextValue 13 4 '4.577'
Set-TextValue 13 5 '  +1.01%  '
Set-TextValue 14 4 '0.6006'
Set-TextValue 14 5 '  -1.51%  '
Set-TextValue 15 4 '77.33'
Set-TextValue 15 5 '  +0.13%  '
Set-TextValue 16 5 '  +0.01%  '
Set-TextValue 17 4 '26.423.75'
Set-TextValue 17 5 '  -0.39%  '
Set-TextValue 18 4 '0.9994'
Set-TextValue 19 4 '0.000007188'
Set-TextValue 19 5 '  +3.52%  '
Set-TextValue 20 4 '11.36'
Set-TextValue 20 5 '  -1.58%  '
Set-TextValue 21 4 '1.940.32'
Set-TextValue 21 5 '  -0.60%  '
Set-TextValue 22 4 '4.490'
Set-TextValue 22 5 '  -0.87%  '
Set-TextValue 23 4 '8.593'
Set-TextValue 23 5 '  -2.38%  '
Set-TextValue 24 4 '5.173'
Set-TextValue 24 5 '  -1.50%  '
Set-TextValue 25 4 '137.78'
Set-TextValue 25 5 '  +0.66%  '
Set-TextValue 26 4 '15.25'
Set-TextValue 26 5 '  -0.60%  '
Set-TextValue 28 4 '107.15'
Set-TextValue 28 5 '  -0.11%  '
Set-TextValue 29 4 '1.717'
Set-TextValue 29 5 '  -3.55%  '
Set-TextValue 30 4 '3.950'
Set-TextValue 30 5 '  -0.73%  '
Set-TextValue 31 4 '0.07977'
Set-TextValue 31 5 '  -0.67%  '
Set-TextValue 32 4 '3.672'
Set-TextValue 32 5 '  -0.69%  '
Set-TextValue 33 5 '  +0.61%  '
Set-TextValue 34 2 'Frax'
Set-TextValue 34 3 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 34 4 '0.9992'
Set-TextValue 34 5 '  -0.01%  '
Set-TextValue 35 2 'HuobiToken'
Set-TextValue 35 3 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 35 4 '2.602'
Set-TextValue 35 5 '  -0.57%  '
Set-TextValue 36 2 'ARBITRUM'
Set-TextValue 36 3 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 36 4 '0.9967'
Set-TextValue 36 5 '  -0.63%  '
Set-TextValue 37 2 'ImmutableX'
Set-TextValue 37 3 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 37 4 '0.6264'
Set-TextValue 37 5 '  +0.04%  '
Set-TextValue 38 2 'TrustWalletToken'
Set-TextValue 38 3 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 38 4 '0.9201'
Set-TextValue 38 5 '  +1.20%  '
Set-TextValue 39 4 '1.964'
Set-TextValue 39 5 '  -5.56%  '
Set-TextValue 40 2 'MXToken'
Set-TextValue 40 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 40 4 '2.393'
Set-TextValue 40 5 '  +0.30%  '
Set-TextValue 41 2 'PaxDollar'
Set-TextValue 41 3 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 41 4 '0.9996'
Set-TextValue 41 5 '  -0.37%  '
Set-TextValue 42 2 'VeChain'
Set-TextValue 42 3 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 42 4 '0.01489'
Set-TextValue 42 5 '  -0.98%  '
Set-TextValue 43 2 'Quant'
Set-TextValue 43 3 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 43 4 '100.02'
Set-TextValue 43 5 '  -2.36%  '
Set-TextValue 44 2 'FraxShare'
Set-TextValue 44 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 44 4 '5.350'
Set-TextValue 44 5 '  -3.25%  '
Set-TextValue 45 2 'TheSandbox'
Set-TextValue 45 3 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 45 4 '0.3849'
Set-TextValue 45 5 '  -0.75%  '
Set-TextValue 46 2 'Aptos'
Set-TextValue 46 3 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 46 4 '6.726'
Set-TextValue 46 5 '  -3.68%  '
Set-TextValue 47 2 'Algorand'
Set-TextValue 47 3 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 47 4 '0.1165'
Set-TextValue 47 5 '  -1.28%  '
Set-TextValue 48 2 'Cronos'
Set-TextValue 48 3 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 48 4 '0.05364'
Set-TextValue 48 5 '  -0.14%  '
Set-TextValue 49 2 'Elrond'
Set-TextValue 49 3 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue 49 4 '30.14'
Set-TextValue 49 5 '  -1.89%  '
Set-TextValue 50 2 'EnergySwap'
Set-TextValue 50 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 50 4 '7.705'
Set-TextValue 50 5 '  -2.02%  '
Set-TextValue 51 2 'NEARProtocol'
Set-TextValue 51 3 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 51 4 '1.237'
Set-TextValue 51 5 '  -1.05%  '
